# environmental-management-thresholds.xlsx edit
#
# Semantic change (per the commit's xml diff):
#   1. Delete data row for option "P1 50" (original row 2) — that option
#      doesn't belong in this sheet's series.
#   2. Delete the duplicate/bad "Y 12,22" row for "March 2022" whose max
#      mark was recorded as 0 (original row 12) — a bad scrape.
#   3. Delete the "subject" column (column L) — redundant, every row in
#      this per-subject sheet is the same value.
#   4. The grade-boundary "g" column (J) values are now stored as text
#      (quote-prefixed numbers) instead of numbers, matching how the
#      scraper now emits the lowest boundary.
#
# Row deletions are done bottom-up first so the row numbers used below
# stay valid while both deletes happen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the "Y 12,22 / March 2022 / max mark 0" row (original row 12).
$ws.Rows(12).Delete()

# 2) Delete the "P1 50" row (original row 2).
$ws.Rows(2).Delete()

# 3) Delete the "subject" column (original column L); after the row
#    deletes above this is still column L (rows don't affect columns).
$ws.Columns(12).Delete()

# 4) Re-type column J (the "g" grade-boundary column) as text for every
#    remaining data row (now rows 2-25), using a quote prefix so Excel
#    stores a literal string instead of re-parsing it as a number.
$gValues = @{
    2  = "25"
    3  = "22"
    4  = "25"
    5  = "20"
    6  = "24"
    7  = "18"
    8  = "22"
    9  = "28"
    10 = "28"
    11 = "24"
    12 = "24"
    13 = "24"
    14 = "23"
    15 = "23"
    16 = "23"
    17 = "25"
    18 = "28"
    19 = "25"
    20 = "37"
    21 = "37"
    22 = "35"
    23 = "23"
    24 = "22"
    25 = "24"
}

foreach ($r in $gValues.Keys) {
    $ws.Cells.Item($r, 10).Value = "'" + $gValues[$r]
}
